$d = $word.ActiveDocument

# --- Paragraph 1: the **ID__...__ID** placeholder line ---
$p = $d.Paragraphs(1)
$r = $p.Range

# The paragraph currently ends in a lone trailing-space run right before the
# paragraph mark (".Range.Text" includes the paragraph mark as its final
# character, so the space sits at End-2 .. End-1). Remove that whole run.
$spaceRange = $d.Range($r.End - 2, $r.End - 1)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Delete()
}

# Update the placeholder identifier text itself.
$d.Content.Find.Execute("**ID__AFFARS_mp_5315_3_topic_2__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_MP_5315_3_1__ID**", 2)

# Re-fetch the paragraph/range (text length changed) and apply the new
# paragraph-level formatting: a thin paragraph border on all four sides and
# a slightly larger left indent.
$p = $d.Paragraphs(1)
$p.Range.ParagraphFormat.LeftIndent = 11.25

$borders = $p.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
